# MicroDustBuildingConfig.xlsx edit:
#  - Rename sheet "UnitProto" -> "Buildings"
#  - Add three new columns (J: Stone, K: Wood, L: Metal) holding the
#    resource cost required to build/upgrade each building, mirroring the
#    MaxLevel (I) column for every row except the Storage building (row 6),
#    which carries its own per-resource amounts.
#  - Update the view (selection / scroll position) to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- rename sheet -----------------------------------------------------
$ws.Name = "Buildings"

# --- new header cells (rows 3-5), copying format from col I -----------
$ws.Range("I3").Copy($ws.Range("J3"))
$ws.Range("I4").Copy($ws.Range("J4"))
$ws.Range("I5").Copy($ws.Range("J5"))
$ws.Range("J3:J5").Copy($ws.Range("K3:L5"))

$ws.Range("J3").Value2 = "Stone"
$ws.Range("K3").Value2 = "Wood"
$ws.Range("L3").Value2 = "Metal"

$ws.Range("J4").Value2 = "Stone"
$ws.Range("K4").Value2 = "Wood"
$ws.Range("L4").Value2 = "Metal"

$ws.Range("J5").Value2 = "int[]"
$ws.Range("K5").Value2 = "int[]"
$ws.Range("L5").Value2 = "int[]"

# --- resource-cost data, rows 6-45 -------------------------------------
# each entry: row, Stone, Wood, Metal
$data = @(
    @(6, "20,20,20", "20,30,20", "20,10"),
    @(7, 30, 30, 30),
    @(8, 30, 30, 30),
    @(9, 30, 30, 30),
    @(10, 30, 30, 30),
    @(11, 20, 20, 20),
    @(12, 3, 3, 3),
    @(13, 4, 4, 4),
    @(14, 5, 5, 5),
    @(15, 5, 5, 5),
    @(16, 4, 4, 4),
    @(17, 5, 5, 5),
    @(18, 5, 5, 5),
    @(19, 10, 10, 10),
    @(20, 10, 10, 10),
    @(21, 10, 10, 10),
    @(22, 10, 10, 10),
    @(23, 10, 10, 10),
    @(24, 5, 5, 5),
    @(25, 10, 10, 10),
    @(26, 10, 10, 10),
    @(27, 10, 10, 10),
    @(28, 10, 10, 10),
    @(29, 20, 20, 20),
    @(30, 3, 3, 3),
    @(31, 10, 10, 10),
    @(32, 10, 10, 10),
    @(33, 10, 10, 10),
    @(34, 10, 10, 10),
    @(35, 10, 10, 10),
    @(36, 5, 5, 5),
    @(37, 5, 5, 5),
    @(38, 10, 10, 10),
    @(39, 3, 3, 3),
    @(40, 10, 10, 10),
    @(41, 5, 5, 5),
    @(42, 4, 4, 4),
    @(43, 5, 5, 5),
    @(44, 5, 5, 5),
    @(45, 5, 5, 5)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $ws.Range("I$r").Copy($ws.Range("J$r"))
    $ws.Range("I$r").Copy($ws.Range("K$r"))
    $ws.Range("I$r").Copy($ws.Range("L$r"))
    $ws.Cells.Item($r, 10).Value2 = $entry[1]
    $ws.Cells.Item($r, 11).Value2 = $entry[2]
    $ws.Cells.Item($r, 12).Value2 = $entry[3]
}

# --- column widths for the new columns (match column I) ---------------
$ws.Columns("J").ColumnWidth = $ws.Columns("I").ColumnWidth
$ws.Columns("K").ColumnWidth = $ws.Columns("I").ColumnWidth
$ws.Columns("L").ColumnWidth = $ws.Columns("I").ColumnWidth

# --- view: selection on L6, no forced scroll ---------------------------
$ws.Range("L6").Select()
